# Update database and change read_price algorithm.
# Row 27 ("سود هر سهم بر اساس آخرین سرمایه" / EPS based on latest capital)
# values are recomputed with the new read_price algorithm (values / 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("D27").Value = 186
$ws.Range("E27").Value = 372
$ws.Range("F27").Value = 502
$ws.Range("G27").Value = 129
$ws.Range("H27").Value = 294
$ws.Range("I27").Value = 510
$ws.Range("J27").Value = 930
$ws.Range("K27").Value = 291
$ws.Range("L27").Value = 483
$ws.Range("M27").Value = 725
